$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K6").Value = 1500075
$ws.Range("M6").Value = -1499963
$ws.Range("H6").Value = 500025
$ws.Range("I6").Value = 500025

$ws.Range("M64").Value = -3671.6667
$ws.Range("I64").Value = 3919.6667
$ws.Range("L64").Value = 7500
$ws.Range("N64").Value = -7996
$ws.Range("H64").Value = 5351.8
$ws.Range("J64").Value = 7500
$ws.Range("K64").Value = 3919.6667

$ws.Range("I67").Value = 3919.6667
$ws.Range("L67").Value = 7500
$ws.Range("K67").Value = 3919.6667
$ws.Range("M67").Value = -3061.6667
$ws.Range("N67").Value = -9216
$ws.Range("H67").Value = 5351.8
$ws.Range("J67").Value = 7500

$ws.Range("I116").Value = 4611.222
$ws.Range("N116").Value = -10302
$ws.Range("K116").Value = 4611.222
$ws.Range("M116").Value = -1169.222
$ws.Range("L116").Value = 3418
$ws.Range("H116").Value = 4244.077
$ws.Range("J116").Value = 3418

$ws.Range("I132").Value = 6232.8
$ws.Range("M132").Value = -16168.4
$ws.Range("K132").Value = 18698.4
$ws.Range("L132").Value = 12088.7145
$ws.Range("H132").Value = 5531.773
$ws.Range("N132").Value = -17148.7145
$ws.Range("J132").Value = 4029.5715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K32").Value = 5833.3335
$ws.Range("M32").Value = -5546.3335
$ws.Range("H32").Value = 6637.8965
$ws.Range("I32").Value = 5833.3335

$ws.Range("H61").Value = 6082.4443
$ws.Range("K61").Value = 6392
$ws.Range("I61").Value = 6392
$ws.Range("M61").Value = -6180

$ws.Range("I74").Value = 2341
$ws.Range("M74").Value = -1467
$ws.Range("K74").Value = 2341
$ws.Range("H74").Value = 2808.4707

$ws.Range("N76").Value = -48342
$ws.Range("H76").Value = 47666
$ws.Range("L76").Value = 47666
$ws.Range("J76").Value = 47666

$ws.Range("H77").Value = 2808.4707
$ws.Range("M77").Value = -7337
$ws.Range("I77").Value = 2341
$ws.Range("K77").Value = 11705

$ws.Range("L79").Value = 47666
$ws.Range("J79").Value = 47666
$ws.Range("H79").Value = 47666
$ws.Range("N79").Value = -50006

$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("J111").Value = 0
$ws.Range("H111").Value = 0

$ws.Range("H114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H136").Value = 6082.4443
$ws.Range("M136").Value = -16626
$ws.Range("I136").Value = 6392
$ws.Range("K136").Value = 19176

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N25").Value = -18470
$ws.Range("H25").Value = 9150
$ws.Range("L25").Value = 18000
$ws.Range("J25").Value = 18000

$ws.Range("K99").Value = 950
$ws.Range("L99").Value = 2998
$ws.Range("M99").Value = 548
$ws.Range("H99").Value = 2315.3333
$ws.Range("J99").Value = 2998
$ws.Range("N99").Value = -5994
$ws.Range("I99").Value = 950

$ws.Range("K102").Value = 10000
$ws.Range("H102").Value = 10000
$ws.Range("I102").Value = 10000
$ws.Range("M102").Value = -6755

$ws.Range("I105").Value = 2549.625
$ws.Range("N105").Value = -6444
$ws.Range("K105").Value = 2549.625
$ws.Range("J105").Value = 2950
$ws.Range("M105").Value = -802.625
$ws.Range("L105").Value = 2950
$ws.Range("H105").Value = 2629.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K99").Value = 5995.2666
$ws.Range("M99").Value = -4497.2666
$ws.Range("H99").Value = 5995.2666
$ws.Range("I99").Value = 5995.2666

$ws.Range("I126").Value = 5995.2666
$ws.Range("M126").Value = -15515.7998
$ws.Range("K126").Value = 17985.7998
$ws.Range("H126").Value = 5995.2666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I5").Value = 885.1667
$ws.Range("K5").Value = 2655.5001
$ws.Range("H5").Value = 885.1667
$ws.Range("M5").Value = -2543.5001

$ws.Range("I92").Value = 320.18182
$ws.Range("M92").Value = 287.45454
$ws.Range("N92").Value = -2643
$ws.Range("H92").Value = 278.46155
$ws.Range("L92").Value = 147
$ws.Range("K92").Value = 960.54546
$ws.Range("J92").Value = 49

$ws.Range("L98").Value = 1497
$ws.Range("N98").Value = -4493
$ws.Range("J98").Value = 499
$ws.Range("H98").Value = 499

$ws.Range("J122").Value = 0
$ws.Range("H122").Value = 1999
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("I132").Value = 999
$ws.Range("M132").Value = -6461
$ws.Range("H132").Value = 999
$ws.Range("K132").Value = 8991

$ws.Range("I135").Value = 885.1667
$ws.Range("K135").Value = 7966.5003
$ws.Range("H135").Value = 885.1667
$ws.Range("M135").Value = -5431.5003

$ws.Range("H140").Value = 1003929.8
$ws.Range("N140").ClearContents()
$ws.Range("M140").Value = -3006609.4
$ws.Range("L140").Value = 0
$ws.Range("I140").Value = 1003929.8
$ws.Range("K140").Value = 3011789.4
$ws.Range("J140").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I132").Value = 2926.1428
$ws.Range("M132").Value = -6248.428400000001
$ws.Range("K132").Value = 8778.428400000001
$ws.Range("L132").Value = 9303
$ws.Range("H132").Value = 2965
$ws.Range("N132").Value = -14363
$ws.Range("J132").Value = 3101

$ws.Range("N134").Value = -116070
$ws.Range("H134").Value = 37000
$ws.Range("J134").Value = 37000
$ws.Range("L134").Value = 111000

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K10").Value = 600.6667
$ws.Range("M10").Value = -460.6667
$ws.Range("H10").Value = 600.3333
$ws.Range("I10").Value = 600.6667

$ws.Range("H22").Value = 1052.1666
$ws.Range("M22").Value = -791
$ws.Range("J22").Value = 984.5
$ws.Range("N22").Value = -1574.5
$ws.Range("I22").Value = 1086
$ws.Range("L22").Value = 984.5
$ws.Range("K22").Value = 1086

$ws.Range("I27").Value = 1086
$ws.Range("H27").Value = 1052.1666
$ws.Range("N27").Value = -1198.5
$ws.Range("J27").Value = 984.5
$ws.Range("K27").Value = 1086
$ws.Range("L27").Value = 984.5
$ws.Range("M27").Value = -979

$ws.Range("K40").Value = 2373
$ws.Range("I40").Value = 2373
$ws.Range("H40").Value = 3029.6667
$ws.Range("M40").Value = -2237

$ws.Range("H68").Value = 4500

$ws.Range("H71").Value = 4500

$ws.Range("H130").Value = 62990
$ws.Range("N130").Value = -73030
$ws.Range("J130").Value = 62990
$ws.Range("L130").Value = 62990

$ws.Range("I132").Value = 4981.8335
$ws.Range("M132").Value = -12415.5005
$ws.Range("N132").Value = -20060
$ws.Range("L132").Value = 15000
$ws.Range("H132").Value = 4984.4287
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 14945.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J10").Value = 16999.5
$ws.Range("L10").Value = 16999.5
$ws.Range("N10").Value = -17337.5
$ws.Range("H10").Value = 16999.5

$ws.Range("I49").Value = 12000
$ws.Range("K49").Value = 12000
$ws.Range("H49").Value = 12000
$ws.Range("M49").Value = -11770

$ws.Range("H62").Value = 9665.666999999999
$ws.Range("J62").Value = 9665.666999999999
$ws.Range("L62").Value = 9665.666999999999
$ws.Range("N62").Value = -10913.667

$ws.Range("H65").Value = 9665.666999999999
$ws.Range("J65").Value = 9665.666999999999
$ws.Range("N65").Value = -54568.335
$ws.Range("L65").Value = 48328.335

$ws.Range("L69").Value = 100271
$ws.Range("H69").Value = 72758.5
$ws.Range("J69").Value = 100271
$ws.Range("N69").Value = -101769

$ws.Range("L72").Value = 300813
$ws.Range("J72").Value = 100271
$ws.Range("N72").Value = -308301
$ws.Range("H72").Value = 72758.5
